# Create a new "mean_ES" summary worksheet at the front of the workbook,
# with the mean of the "ES" column (col H) pulled from each of the
# per-form worksheets, for every scale row.

$wb = $excel.ActiveWorkbook

# Snapshot the existing (pre-insert) worksheets in tab order before we add
# the new summary sheet, so we know which sheets to aggregate over.
$sourceSheets = @()
foreach ($s in $wb.Worksheets) {
    $sourceSheets += $s
}

# Scale names, in row order (from column A of each source sheet).
$scales = @("TOT_raw", "SOC_raw", "VIS_raw", "HEA_raw", "TOU_raw", "TS_raw", "BOD_raw", "BAL_raw", "PLA_raw")

# Compute the across-sheet mean of column H (ES) for each scale row FIRST,
# while $sourceSheets still refers to the original (pre-insert) sheets --
# worksheet references in this runtime resolve positionally, so they would
# shift underneath us once a new sheet is inserted at the front.
$meanValues = @()
for ($i = 0; $i -lt $scales.Length; $i++) {
    $row = $i + 2

    $total = 0.0
    $count = 0
    foreach ($src in $sourceSheets) {
        $total = $total + $src.Cells.Item($row, 8).Value()
        $count = $count + 1
    }
    $meanValues += [Math]::Round($total / $count, 3)
}

# Now insert the new sheet at the very front of the workbook.
$meanSheet = $wb.Worksheets.Add()
$meanSheet.Name = "mean_ES"

# Header row.
$meanSheet.Range("A1").Value = "scale"
$meanSheet.Range("B1").Value = "mean_ES"
$meanSheet.Range("A1:B1").Font.Bold = $true
$meanSheet.Range("A1:B1").HorizontalAlignment = -4108

# One row per scale.
for ($i = 0; $i -lt $scales.Length; $i++) {
    $row = $i + 2
    $meanSheet.Range("A" + $row).Value = $scales[$i]
    $meanSheet.Range("B" + $row).Value = $meanValues[$i]
}

$wb.Worksheets.Item("mean_ES").Activate()
